$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The B/C rows within each year group were swapped (e.g. "2016年B" <-> "2016年C"),
# carrying their A:E data with them. Do this for each affected pair of rows.
$swapPairs = @(3,4), (7,8), (11,12), (15,16)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    for ($col = 1; $col -le 5; $col++) {
        $c1 = $ws.Cells.Item($r1, $col)
        $c2 = $ws.Cells.Item($r2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value2 = $v2
        $c2.Value2 = $v1
    }
}

# Remove the now-unwanted "产销率" and "销售量" columns (F and G), header + data.
$ws.Range("F1:G17").Delete()
